$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 54.2
$ws.Range("E3").Value = 56.9
$ws.Range("E4").Value = 43.6
$ws.Range("E5").Value = 41.2
$ws.Range("E6").Value = 46.2

$ws.Range("G2").Value = 50
$ws.Range("G3").Value = 40
$ws.Range("G4").Value = 30
$ws.Range("G5").Value = 20
$ws.Range("G6").Value = 30

$ws.Range("K2").Value = 67.2
$ws.Range("K3").Value = 57.8
$ws.Range("K4").Value = 56
$ws.Range("K5").Value = 55.8
$ws.Range("K6").Value = 52

$ws.Range("N2").Value = 85.87127175646313
$ws.Range("N3").Value = 85.87127175646313
$ws.Range("N4").Value = 85.87127175646313
$ws.Range("N5").Value = 85.87127175646313
$ws.Range("N6").Value = 85.87127175646313
